$d = $word.ActiveDocument

# This document's body XML, as it exists right now (pristine, byte-for-byte).
# We apply three literal text substitutions to it -- one per hunk of the target
# diff -- and then push the whole body back in a single Range.InsertXML call on
# Document.Content. A full-content InsertXML replaces the body in place instead of
# inserting a sibling paragraph, which is what happens if you try to InsertXML a
# Range that sits on the document's very last (always-collapsed) paragraph alone.
$orig = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<w:document xmlns:wpc="http://schemas.microsoft.com/office/word/2010/wordprocessingCanvas" xmlns:mc="http://schemas.openxmlformats.org/markup-compatibility/2006" xmlns:o="urn:schemas-microsoft-com:office:office" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math" xmlns:v="urn:schemas-microsoft-com:vml" xmlns:wp14="http://schemas.microsoft.com/office/word/2010/wordprocessingDrawing" xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" xmlns:w10="urn:schemas-microsoft-com:office:word" xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" xmlns:w15="http://schemas.microsoft.com/office/word/2012/wordml" xmlns:wpg="http://schemas.microsoft.com/office/word/2010/wordprocessingGroup" xmlns:wpi="http://schemas.microsoft.com/office/word/2010/wordprocessingInk" xmlns:wne="http://schemas.microsoft.com/office/word/2006/wordml" xmlns:wps="http://schemas.microsoft.com/office/word/2010/wordprocessingShape" mc:Ignorable="w14 w15 wp14"><w:body><w:p w:rsidR="001D2FA1" w:rsidRDefault="00FD2D63"><w:r><w:t>Exercises:</w:t></w:r></w:p><w:p w:rsidR="00FD2D63" w:rsidRDefault="00FD2D63" w:rsidP="00FD2D63"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Declare 2 strings</w:t></w:r></w:p><w:p w:rsidR="00FD2D63" w:rsidRDefault="00FD2D63" w:rsidP="004D6BED"><w:pPr><w:pStyle w:val="ListParagraph"/></w:pPr><w:proofErr w:type="gramStart"/><w:r><w:t>string</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> name1 = “Jason” , name2 = “Curry;</w:t></w:r></w:p><w:p w:rsidR="00FD2D63" w:rsidRDefault="00FD2D63" w:rsidP="00FD2D63"><w:pPr><w:pStyle w:val="ListParagraph"/></w:pPr><w:proofErr w:type="gramStart"/><w:r><w:t>write</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> the C++ code to display the output with this 2 strings.</w:t></w:r></w:p><w:p w:rsidR="00FD2D63" w:rsidRDefault="00FD2D63" w:rsidP="004D6BED"><w:pPr><w:ind w:firstLine="720"/></w:pPr><w:r><w:t>Output:</w:t></w:r></w:p><w:p w:rsidR="001E42C0" w:rsidRDefault="001E42C0" w:rsidP="00FD2D63"><w:r><w:tab/></w:r><w:r w:rsidR="004D6BED"><w:tab/></w:r><w:r><w:t>Jason</w:t></w:r><w:r w:rsidR="00FD2D63"><w:t>******************************Curry</w:t></w:r></w:p><w:p w:rsidR="001E42C0" w:rsidRDefault="001E42C0" w:rsidP="004D6BED"><w:pPr><w:ind w:firstLine="720"/></w:pPr><w:r><w:t xml:space="preserve">Requirement: </w:t></w:r></w:p><w:p w:rsidR="001E42C0" w:rsidRDefault="001E42C0" w:rsidP="001E42C0"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Using </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>setw</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">() function to make a 20 digit length space, then set name1 in the space and align to left </w:t></w:r></w:p><w:p w:rsidR="001E42C0" w:rsidRDefault="001E42C0" w:rsidP="001E42C0"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Do the same thing to name2, but set </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>the align</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> to right.</w:t></w:r></w:p><w:p w:rsidR="001E42C0" w:rsidRDefault="001E42C0" w:rsidP="001E42C0"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Using </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>setfill</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>() function to add  all the * between name1 and name2</w:t></w:r></w:p><w:p w:rsidR="008E5F0E" w:rsidRDefault="008E5F0E" w:rsidP="008E5F0E"/><w:p w:rsidR="00954FA7" w:rsidRDefault="00954FA7" w:rsidP="008E5F0E"/><w:p w:rsidR="008E5F0E" w:rsidRDefault="008E5F0E" w:rsidP="008E5F0E"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Prompt user to input his full name and just read the last name </w:t></w:r></w:p><w:p w:rsidR="00B17A7F" w:rsidRDefault="00B17A7F" w:rsidP="00B17A7F"><w:pPr><w:pStyle w:val="ListParagraph"/></w:pPr><w:r><w:t>Example:</w:t></w:r></w:p><w:p w:rsidR="00B17A7F" w:rsidRDefault="00B17A7F" w:rsidP="00B17A7F"><w:pPr><w:pStyle w:val="ListParagraph"/></w:pPr><w:r><w:tab/><w:t>Assume the input is: Michael Jordan</w:t></w:r></w:p><w:p w:rsidR="00B17A7F" w:rsidRDefault="00B17A7F" w:rsidP="00B17A7F"><w:pPr><w:pStyle w:val="ListParagraph"/></w:pPr><w:r><w:tab/><w:t xml:space="preserve">Create a string </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>last_name</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>, and read the last name only from input.</w:t></w:r></w:p><w:p w:rsidR="00B17A7F" w:rsidRDefault="00B17A7F" w:rsidP="00B17A7F"><w:pPr><w:pStyle w:val="ListParagraph"/></w:pPr><w:r><w:tab/><w:t>(</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>using</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>cin.ignore</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>())</w:t></w:r></w:p><w:p w:rsidR="00DC09AC" w:rsidRDefault="00DC09AC" w:rsidP="00DC09AC"/><w:p w:rsidR="00954FA7" w:rsidRDefault="00954FA7" w:rsidP="00DC09AC"><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p><w:p w:rsidR="003B6A50" w:rsidRDefault="00516BCB" w:rsidP="003B6A50"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Show all the character of asci code from 1 ~ 256</w:t></w:r><w:r w:rsidR="003B6A50"><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidR="00557FCF"><w:t xml:space="preserve"> (do later after for loop reviewed)</w:t></w:r></w:p><w:p w:rsidR="003B6A50" w:rsidRDefault="003B6A50" w:rsidP="00A87E54"><w:pPr><w:pStyle w:val="ListParagraph"/></w:pPr><w:r><w:t>Sample output:</w:t></w:r></w:p><w:p w:rsidR="003B6A50" w:rsidRDefault="003B6A50" w:rsidP="003B6A50"><w:pPr><w:pStyle w:val="ListParagraph"/></w:pPr><w:r><w:rPr><w:noProof/></w:rPr><w:drawing><wp:inline distT="0" distB="0" distL="0" distR="0"><wp:extent cx="2542857" cy="923810"/><wp:effectExtent l="0" t="0" r="0" b="0"/><wp:docPr id="1" name="Picture 1"/><wp:cNvGraphicFramePr><a:graphicFrameLocks xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" noChangeAspect="1"/></wp:cNvGraphicFramePr><a:graphic xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main"><a:graphicData uri="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:pic xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:nvPicPr><pic:cNvPr id="1" name="test.png"/><pic:cNvPicPr/></pic:nvPicPr><pic:blipFill><a:blip r:embed="rId5"><a:extLst><a:ext uri="{28A0092B-C50C-407E-A947-70E740481C1C}"><a14:useLocalDpi xmlns:a14="http://schemas.microsoft.com/office/drawing/2010/main" val="0"/></a:ext></a:extLst></a:blip><a:stretch><a:fillRect/></a:stretch></pic:blipFill><pic:spPr><a:xfrm><a:off x="0" y="0"/><a:ext cx="2542857" cy="923810"/></a:xfrm><a:prstGeom prst="rect"><a:avLst/></a:prstGeom></pic:spPr></pic:pic></a:graphicData></a:graphic></wp:inline></w:drawing></w:r></w:p><w:p w:rsidR="00A87E54" w:rsidRDefault="00A87E54" w:rsidP="00A87E54"/><w:p w:rsidR="00A87E54" w:rsidRDefault="00A87E54" w:rsidP="00852BAD"><w:pPr><w:pStyle w:val="ListParagraph"/></w:pPr></w:p><w:p w:rsidR="004D6BED" w:rsidRDefault="004D6BED" w:rsidP="004D6BED"/><w:p w:rsidR="004D6BED" w:rsidRDefault="004D6BED" w:rsidP="004D6BED"/><w:p w:rsidR="00FD2D63" w:rsidRDefault="00FD2D63" w:rsidP="001E42C0"><w:pPr><w:pStyle w:val="ListParagraph"/><w:ind w:left="1080"/><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr></w:pPr></w:p><w:sectPr w:rsidR="00FD2D63"><w:pgSz w:w="12240" w:h="15840"/><w:pgMar w:top="1440" w:right="1440" w:bottom="1440" w:left="1440" w:header="720" w:footer="720" w:gutter="0"/><w:cols w:space="720"/><w:docGrid w:linePitch="360"/></w:sectPr></w:body></w:document>'

# --- Hunk 2: drop the old _GoBack bookmark that sits alone in its own paragraph.
$oldBookmarkPara = '<w:p w:rsidR="00954FA7" w:rsidRDefault="00954FA7" w:rsidP="00DC09AC"><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'
$newBookmarkPara = '<w:p w:rsidR="00954FA7" w:rsidRDefault="00954FA7" w:rsidP="00DC09AC"/>'
if ($orig.IndexOf($oldBookmarkPara) -lt 0) { throw "old bookmark paragraph not found" }
$updated = $orig.Replace($oldBookmarkPara, $newBookmarkPara)

# --- Hunk 1: "Exercises:" -> "Exercise" + a fresh _GoBack bookmark + ":".
$oldTitleRun = '<w:r><w:t>Exercises:</w:t></w:r>'
$newTitleRun = '<w:r><w:t>Exercise</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:t>:</w:t></w:r>'
if ($updated.IndexOf($oldTitleRun) -lt 0) { throw "title run not found" }
$updated = $updated.Replace($oldTitleRun, $newTitleRun)

# --- Hunk 3: drop the stray eastAsia font hint on the trailing empty paragraph.
$oldRpr = '<w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr>'
if ($updated.IndexOf($oldRpr) -lt 0) { throw "rFonts hint not found" }
$updated = $updated.Replace($oldRpr, '')

$pkg = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData>' + $updated + '</pkg:xmlData></pkg:part></pkg:package>'
$d.Content.InsertXML($pkg)

Write-Output "paragraphs=$($d.Paragraphs.Count)"
